# Scheduled data refresh: update market-price-derived columns (H-N) across all
# job sheets in the Chocobo Profits workbook. Values sourced from the latest
# Universalis market snapshot; only numeric columns H..N are touched, A..G are
# static leve metadata and are left untouched.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(98, 8).Value = 2731.2903  # H98: 2626.3225 -> 2731.2903
$ws.Cells.Item(98, 9).Value = 1487.7778  # I98: 1323.6923 -> 1487.7778
$ws.Cells.Item(98, 10).Value = 11125  # J98: 9400 -> 11125
$ws.Cells.Item(98, 11).Value = 1487.7778  # K98: 1323.6923 -> 1487.7778
$ws.Cells.Item(98, 12).Value = 11125  # L98: 9400 -> 11125
$ws.Cells.Item(98, 13).Value = 10.22219999999993  # M98: 174.3077000000001 -> 10.22219999999993
$ws.Cells.Item(98, 14).Value = -14121  # N98: -12396 -> -14121

$ws.Cells.Item(112, 8).Value = 1632.6548  # H112: 1564.8877 -> 1632.6548
$ws.Cells.Item(112, 10).Value = 1643.8916  # J112: 1573.8041 -> 1643.8916
$ws.Cells.Item(112, 12).Value = 4931.6748  # L112: 4721.4123 -> 4931.6748
$ws.Cells.Item(112, 14).Value = -7147.6748  # N112: -6937.4123 -> -7147.6748

$ws.Cells.Item(122, 8).Value = 2731.2903  # H122: 2626.3225 -> 2731.2903
$ws.Cells.Item(122, 9).Value = 1487.7778  # I122: 1323.6923 -> 1487.7778
$ws.Cells.Item(122, 10).Value = 11125  # J122: 9400 -> 11125
$ws.Cells.Item(122, 11).Value = 4463.3334  # K122: 3971.0769 -> 4463.3334
$ws.Cells.Item(122, 12).Value = 33375  # L122: 28200 -> 33375
$ws.Cells.Item(122, 13).Value = -2013.3334  # M122: -1521.0769 -> -2013.3334
$ws.Cells.Item(122, 14).Value = -38275  # N122: -33100 -> -38275

$ws.Cells.Item(137, 8).Value = 2837.457  # H137: 2561.9211 -> 2837.457
$ws.Cells.Item(137, 9).Value = 1840.3928  # I137: 1740.4517 -> 1840.3928
$ws.Cells.Item(137, 10).Value = 6825.7144  # J137: 6199.857 -> 6825.7144
$ws.Cells.Item(137, 11).Value = 5521.178400000001  # K137: 5221.355100000001 -> 5521.178400000001
$ws.Cells.Item(137, 12).Value = 20477.1432  # L137: 18599.571 -> 20477.1432
$ws.Cells.Item(137, 13).Value = -2971.178400000001  # M137: -2671.355100000001 -> -2971.178400000001
$ws.Cells.Item(137, 14).Value = -25577.1432  # N137: -23699.571 -> -25577.1432

$ws.Cells.Item(138, 8).Value = 2753.81  # H138: 2780.5251 -> 2753.81
$ws.Cells.Item(138, 9).Value = 1470.6364  # I138: 1302.3846 -> 1470.6364
$ws.Cells.Item(138, 10).Value = 2912.4045  # J138: 3003.965 -> 2912.4045
$ws.Cells.Item(138, 11).Value = 4411.9092  # K138: 3907.1538 -> 4411.9092
$ws.Cells.Item(138, 12).Value = 8737.2135  # L138: 9011.895 -> 8737.2135
$ws.Cells.Item(138, 13).Value = 728.0907999999999  # M138: 1232.8462 -> 728.0907999999999
$ws.Cells.Item(138, 14).Value = -19017.2135  # N138: -19291.895 -> -19017.2135

$ws.Cells.Item(139, 8).Value = 41365  # H139: 41375.715 -> 41365
$ws.Cells.Item(139, 10).Value = 41365  # J139: 41375.715 -> 41365
$ws.Cells.Item(139, 12).Value = 41365  # L139: 41375.715 -> 41365
$ws.Cells.Item(139, 14).Value = -51645  # N139: -51655.715 -> -51645

$ws.Cells.Item(141, 8).Value = 1786.4688  # H141: 1907.875 -> 1786.4688
$ws.Cells.Item(141, 9).Value = 1668.7333  # I141: 1725.7587 -> 1668.7333
$ws.Cells.Item(141, 10).Value = 3552.5  # J141: 3668.3333 -> 3552.5
$ws.Cells.Item(141, 11).Value = 5006.199900000001  # K141: 5177.2761 -> 5006.199900000001
$ws.Cells.Item(141, 12).Value = 10657.5  # L141: 11004.9999 -> 10657.5
$ws.Cells.Item(141, 13).Value = 173.8000999999995  # M141: 2.723899999999958 -> 173.8000999999995
$ws.Cells.Item(141, 14).Value = -21017.5  # N141: -21364.9999 -> -21017.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(4, 8).Value = 243.66667  # H4: 500 -> 243.66667
$ws.Cells.Item(4, 9).Value = 265.5  # I4: 500 -> 265.5
$ws.Cells.Item(4, 10).Value = 200  # J4: 0 -> 200
$ws.Cells.Item(4, 11).Value = 265.5  # K4: 500 -> 265.5
$ws.Cells.Item(4, 12).Value = 200  # L4: 0 -> 200
$ws.Cells.Item(4, 13).Value = -149.5  # M4: -384 -> -149.5
$ws.Cells.Item(4, 14).Value = -432  # N4: None -> -432

$ws.Cells.Item(61, 8).Value = 2304.318  # H61: 1998.6 -> 2304.318
$ws.Cells.Item(61, 9).Value = 1766.9166  # I61: 1487.6842 -> 1766.9166
$ws.Cells.Item(61, 10).Value = 2949.2  # J61: 2881.0908 -> 2949.2
$ws.Cells.Item(61, 11).Value = 1766.9166  # K61: 1487.6842 -> 1766.9166
$ws.Cells.Item(61, 12).Value = 2949.2  # L61: 2881.0908 -> 2949.2
$ws.Cells.Item(61, 13).Value = -1554.9166  # M61: -1275.6842 -> -1554.9166
$ws.Cells.Item(61, 14).Value = -3373.2  # N61: -3305.0908 -> -3373.2

$ws.Cells.Item(97, 8).Value = 536.75  # H97: 453.0606 -> 536.75
$ws.Cells.Item(97, 9).Value = 570.8929000000001  # I97: 465.33334 -> 570.8929000000001
$ws.Cells.Item(97, 10).Value = 297.75  # J97: 330.33334 -> 297.75
$ws.Cells.Item(97, 11).Value = 570.8929000000001  # K97: 465.33334 -> 570.8929000000001
$ws.Cells.Item(97, 12).Value = 297.75  # L97: 330.33334 -> 297.75
$ws.Cells.Item(97, 13).Value = -74.89290000000005  # M97: 30.66665999999998 -> -74.89290000000005
$ws.Cells.Item(97, 14).Value = -1289.75  # N97: -1322.33334 -> -1289.75

$ws.Cells.Item(110, 8).Value = 708.8484999999999  # H110: 699.74286 -> 708.8484999999999
$ws.Cells.Item(110, 9).Value = 672.1539  # I110: 663.3929000000001 -> 672.1539
$ws.Cells.Item(110, 11).Value = 672.1539  # K110: 663.3929000000001 -> 672.1539
$ws.Cells.Item(110, 13).Value = 1372.8461  # M110: 1381.6071 -> 1372.8461

$ws.Cells.Item(136, 8).Value = 2304.318  # H136: 1998.6 -> 2304.318
$ws.Cells.Item(136, 9).Value = 1766.9166  # I136: 1487.6842 -> 1766.9166
$ws.Cells.Item(136, 10).Value = 2949.2  # J136: 2881.0908 -> 2949.2
$ws.Cells.Item(136, 11).Value = 5300.7498  # K136: 4463.0526 -> 5300.7498
$ws.Cells.Item(136, 12).Value = 8847.599999999999  # L136: 8643.2724 -> 8847.599999999999
$ws.Cells.Item(136, 13).Value = -2750.7498  # M136: -1913.0526 -> -2750.7498
$ws.Cells.Item(136, 14).Value = -13947.6  # N136: -13743.2724 -> -13947.6

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(43, 8).Value = 89800  # H43: 79800 -> 89800
$ws.Cells.Item(43, 10).Value = 89800  # J43: 79800 -> 89800
$ws.Cells.Item(43, 12).Value = 89800  # L43: 79800 -> 89800
$ws.Cells.Item(43, 14).Value = -90162  # N43: -80162 -> -90162

$ws.Cells.Item(59, 8).Value = 59999  # H59: 0 -> 59999
$ws.Cells.Item(59, 10).Value = 59999  # J59: 0 -> 59999
$ws.Cells.Item(59, 12).Value = 59999  # L59: 0 -> 59999
$ws.Cells.Item(59, 14).Value = -61693  # N59: None -> -61693

$ws.Cells.Item(94, 8).Value = 946.75  # H94: 873.57574 -> 946.75
$ws.Cells.Item(94, 9).Value = 1063.1052  # I94: 892.72 -> 1063.1052
$ws.Cells.Item(94, 10).Value = 701.1111  # J94: 813.75 -> 701.1111
$ws.Cells.Item(94, 11).Value = 1063.1052  # K94: 892.72 -> 1063.1052
$ws.Cells.Item(94, 12).Value = 701.1111  # L94: 813.75 -> 701.1111
$ws.Cells.Item(94, 13).Value = -612.1052  # M94: -441.72 -> -612.1052
$ws.Cells.Item(94, 14).Value = -1603.1111  # N94: -1715.75 -> -1603.1111

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(105, 8).Value = 3148.625  # H105: 1556.95 -> 3148.625
$ws.Cells.Item(105, 9).Value = 4030  # I105: 1076.6 -> 4030
$ws.Cells.Item(105, 10).Value = 2619.8  # J105: 2998 -> 2619.8
$ws.Cells.Item(105, 11).Value = 4030  # K105: 1076.6 -> 4030
$ws.Cells.Item(105, 12).Value = 2619.8  # L105: 2998 -> 2619.8
$ws.Cells.Item(105, 13).Value = -2283  # M105: 670.4000000000001 -> -2283
$ws.Cells.Item(105, 14).Value = -6113.8  # N105: -6492 -> -6113.8

$ws.Cells.Item(134, 8).Value = 8364.235000000001  # H134: 6711.5454 -> 8364.235000000001
$ws.Cells.Item(134, 9).Value = 8784  # I134: 6862 -> 8784
$ws.Cells.Item(134, 10).Value = 7000  # J134: 6200 -> 7000
$ws.Cells.Item(134, 11).Value = 26352  # K134: 20586 -> 26352
$ws.Cells.Item(134, 12).Value = 21000  # L134: 18600 -> 21000
$ws.Cells.Item(134, 13).Value = -23817  # M134: -18051 -> -23817
$ws.Cells.Item(134, 14).Value = -26070  # N134: -23670 -> -26070

$ws.Cells.Item(137, 8).Value = 34674.547  # H137: 36220 -> 34674.547
$ws.Cells.Item(137, 10).Value = 34674.547  # J137: 36220 -> 34674.547
$ws.Cells.Item(137, 12).Value = 34674.547  # L137: 36220 -> 34674.547
$ws.Cells.Item(137, 14).Value = -44874.547  # N137: -46420 -> -44874.547

$ws.Cells.Item(138, 8).Value = 29682.5  # H138: 30080 -> 29682.5
$ws.Cells.Item(138, 10).Value = 29682.5  # J138: 30080 -> 29682.5
$ws.Cells.Item(138, 12).Value = 29682.5  # L138: 30080 -> 29682.5
$ws.Cells.Item(138, 14).Value = -39962.5  # N138: -40360 -> -39962.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(99, 8).Value = 730  # H99: 2064.2856 -> 730
$ws.Cells.Item(99, 9).Value = 730  # I99: 983.3333 -> 730
$ws.Cells.Item(99, 10).Value = 0  # J99: 2875 -> 0
$ws.Cells.Item(99, 11).Value = 2190  # K99: 2949.9999 -> 2190
$ws.Cells.Item(99, 12).Value = 0  # L99: 8625 -> 0
$ws.Cells.Item(99, 13).Value = 56  # M99: -703.9998999999998 -> 56
$ws.Cells.Item(99, 14).ClearContents()  # N99: -13117 -> (removed)

$ws.Cells.Item(114, 8).Value = 3017.3157  # H114: 3658.6667 -> 3017.3157
$ws.Cells.Item(114, 9).Value = 339.1111  # I114: 160 -> 339.1111
$ws.Cells.Item(114, 10).Value = 5427.7  # J114: 7657.143 -> 5427.7
$ws.Cells.Item(114, 11).Value = 1017.3333  # K114: 480 -> 1017.3333
$ws.Cells.Item(114, 12).Value = 16283.1  # L114: 22971.429 -> 16283.1
$ws.Cells.Item(114, 13).Value = 2236.6667  # M114: 2774 -> 2236.6667
$ws.Cells.Item(114, 14).Value = -22791.1  # N114: -29479.429 -> -22791.1

$ws.Cells.Item(126, 8).Value = 2771.875  # H126: 1274.75 -> 2771.875
$ws.Cells.Item(126, 9).Value = 1000  # I126: 0 -> 1000
$ws.Cells.Item(126, 10).Value = 2890  # J126: 1274.75 -> 2890
$ws.Cells.Item(126, 11).Value = 3000  # K126: 0 -> 3000
$ws.Cells.Item(126, 12).Value = 8670  # L126: 3824.25 -> 8670
$ws.Cells.Item(126, 13).Value = 1940  # M126: None -> 1940
$ws.Cells.Item(126, 14).Value = -18550  # N126: -13704.25 -> -18550

$ws.Cells.Item(131, 8).Value = 1119.7755  # H131: 1012.36 -> 1119.7755
$ws.Cells.Item(131, 10).Value = 969.3111  # J131: 910.94366 -> 969.3111
$ws.Cells.Item(131, 12).Value = 2907.9333  # L131: 2732.83098 -> 2907.9333
$ws.Cells.Item(131, 14).Value = -12987.9333  # N131: -12812.83098 -> -12987.9333

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(97, 8).Value = 771.1111  # H97: 709.1111 -> 771.1111
$ws.Cells.Item(97, 9).Value = 743.7692  # I97: 633.17645 -> 743.7692
$ws.Cells.Item(97, 10).Value = 842.2  # J97: 2000 -> 842.2
$ws.Cells.Item(97, 11).Value = 743.7692  # K97: 633.17645 -> 743.7692
$ws.Cells.Item(97, 12).Value = 842.2  # L97: 2000 -> 842.2
$ws.Cells.Item(97, 13).Value = -247.7692  # M97: -137.17645 -> -247.7692
$ws.Cells.Item(97, 14).Value = -1834.2  # N97: -2992 -> -1834.2

$ws.Cells.Item(132, 8).Value = 5218.7  # H132: 3014.1155 -> 5218.7
$ws.Cells.Item(132, 9).Value = 3000  # I132: 1804.7059 -> 3000
$ws.Cells.Item(132, 10).Value = 5773.375  # J132: 5298.5557 -> 5773.375
$ws.Cells.Item(132, 11).Value = 9000  # K132: 5414.1177 -> 9000
$ws.Cells.Item(132, 12).Value = 17320.125  # L132: 15895.6671 -> 17320.125
$ws.Cells.Item(132, 13).Value = -6470  # M132: -2884.1177 -> -6470
$ws.Cells.Item(132, 14).Value = -22380.125  # N132: -20955.6671 -> -22380.125

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(55, 8).Value = 292.8889  # H55: 249.69565 -> 292.8889
$ws.Cells.Item(55, 9).Value = 217.09091  # I55: 187.93333 -> 217.09091
$ws.Cells.Item(55, 10).Value = 412  # J55: 365.5 -> 412
$ws.Cells.Item(55, 11).Value = 217.09091  # K55: 187.93333 -> 217.09091
$ws.Cells.Item(55, 12).Value = 412  # L55: 365.5 -> 412
$ws.Cells.Item(55, 13).Value = -44.09091000000001  # M55: -14.93333000000001 -> -44.09091000000001
$ws.Cells.Item(55, 14).Value = -758  # N55: -711.5 -> -758

$ws.Cells.Item(93, 8).Value = 2159.7646  # H93: 1927.9546 -> 2159.7646
$ws.Cells.Item(93, 9).Value = 1279.1111  # I93: 1131.6154 -> 1279.1111
$ws.Cells.Item(93, 10).Value = 3150.5  # J93: 3078.2222 -> 3150.5
$ws.Cells.Item(93, 11).Value = 1279.1111  # K93: 1131.6154 -> 1279.1111
$ws.Cells.Item(93, 12).Value = 3150.5  # L93: 3078.2222 -> 3150.5
$ws.Cells.Item(93, 13).Value = -31.11110000000008  # M93: 116.3846000000001 -> -31.11110000000008
$ws.Cells.Item(93, 14).Value = -5646.5  # N93: -5574.2222 -> -5646.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(132, 8).Value = 9805879  # H132: 12822870 -> 9805879
$ws.Cells.Item(132, 9).Value = 715.5  # I132: 862.36365 -> 715.5
$ws.Cells.Item(132, 10).Value = 20836688  # J132: 22225676 -> 20836688
$ws.Cells.Item(132, 11).Value = 2146.5  # K132: 2587.09095 -> 2146.5
$ws.Cells.Item(132, 12).Value = 62510064  # L132: 66677028 -> 62510064
$ws.Cells.Item(132, 13).Value = 383.5  # M132: -57.09094999999979 -> 383.5
$ws.Cells.Item(132, 14).Value = -62515124  # N132: -66682088 -> -62515124
